# Update cryptocurrency price/volume data per the Mon May  6 19:25:34 UTC 2024 refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.907.56"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -1.90%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.056.92"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -2.94%  "

# Row 4
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.23%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "588.15"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.68%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "152.17"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +3.83%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.13%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.546"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +2.88%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "3.062.10"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -2.60%  "

# Row 10
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -4.42%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.81"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -1.24%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.459"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.23%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000240"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -3.28%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "37.03"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.47%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.566.09"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -2.86%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.19"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -1.07%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "63.088.82"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -1.34%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.061.95"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -2.66%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "473.07"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +1.01%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.62"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +1.68%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.715"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -2.45%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.52"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.04%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.37"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +1.82%  "

# Row 25
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.64%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "81.13"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.33%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.998"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.24%  "

# Row 28
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +2.12%  "

# Row 29
$ws.Range("B29").Value = "NEARProtocol"
$ws.Range("C29").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.30"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -1.25%  "

# Row 30
$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.67"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -1.65%  "

# Row 31
$ws.Range("B31").Value = "FirstDigitalUSD"
$ws.Range("C31").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.998"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.36%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.19"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -2.57%  "

# Row 33
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +2.74%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "27.18"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -1.91%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0₃0846"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +1.02%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.05"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -2.24%  "

# Row 37
$ws.Range("B37").Value = "Filecoin"
$ws.Range("C37").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.09"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -1.32%  "

# Row 38
$ws.Range("B38").Value = "dogwifhat"
$ws.Range("C38").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.34"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +2.80%  "

# Row 39
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -4.68%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "9.27"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.79%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "50.45"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -1.87%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "442.56"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -4.17%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.284"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -2.94%  "

# Row 44
$ws.Range("B44").Value = "Arweave"
$ws.Range("C44").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "40.21"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +1.49%  "

# Row 45
$ws.Range("B45").Value = "VeChain"
$ws.Range("C45").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0361"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -3.03%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.110"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +1.72%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.796.02"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -4.57%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "131.01"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.65%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "25.10"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +3.86%  "

# Row 51
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.05%  "
